# Apply the commit's edit:
#  1) On the "readme" sheet, reorder the JobNo/Author/Date columns (C/D/E)
#     of Table1 to Date/JobNo/Author, carrying each column's values along.
#  2) On the "Project Information" sheet, refresh the "Date of Analysis"
#     timestamp value.

$wb = $excel.ActiveWorkbook

# --- 1) readme sheet: reorder columns C (JobNo) / D (Author) / E (Date) ---
$readme = $wb.Worksheets.Item("readme")

# Update header row
$readme.Range("C1").Value = "Date"
$readme.Range("D1").Value = "JobNo"
$readme.Range("E1").Value = "Author"

# Rename the table's ListColumns to match the new header order
$table1 = $readme.ListObjects.Item("Table1")
$table1.ListColumns.Item(3).Name = "Date"
$table1.ListColumns.Item(4).Name = "JobNo"
$table1.ListColumns.Item(5).Name = "Author"

# Update each data row: C=Date value, D=JobNo value, E=Author value.
# Column C's new value ("20220325") looks numeric, so Excel would
# auto-convert it to a number under the cell's (General) format. Force it
# to stay TEXT, matching the source data, by briefly switching the cell to
# Text format while assigning, then copy the original cell's formatting
# back over itself (via a same-format sibling cell in column D/E) so the
# cell keeps its original style index instead of acquiring a new one.
$lastRow = $readme.Cells.Item($readme.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cellC = $readme.Cells.Item($r, 3)
    $formatDonor = $readme.Cells.Item($r, 5)  # same style as C, untouched so far

    $cellC.NumberFormat = "@"
    $cellC.Value = "20220325"

    $formatDonor.Copy()
    $cellC.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    $readme.Cells.Item($r, 4).Value = "/c/e"
    $readme.Cells.Item($r, 5).Value = "jovyan"
}

# --- 2) Project Information sheet: update "Date of Analysis" value ---
$projInfo = $wb.Worksheets.Item("Project Information")
$projInfo.Range("B12").Value = "2022-03-25 19:32:57.480416"
